# Pierre Robin Sequence.xlsx edit script
# 1. Refresh the "time_taken" timestamps on the "data" sheet (column F, rows 2-52)
# 2. Add a new "metadata" sheet (after "data") describing the PanelApp query that produced the data

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- 1. Update time_taken column on "data" sheet ---------------------------------------------
$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52)
$timestamps = @("2021-10-05 14:35:15.049590","2021-10-05 14:35:15.049598","2021-10-05 14:35:15.049602","2021-10-05 14:35:15.049604","2021-10-05 14:35:15.049607","2021-10-05 14:35:15.049610","2021-10-05 14:35:15.049612","2021-10-05 14:35:15.049615","2021-10-05 14:35:15.049617","2021-10-05 14:35:15.049620","2021-10-05 14:35:15.049622","2021-10-05 14:35:15.049625","2021-10-05 14:35:15.049627","2021-10-05 14:35:15.049630","2021-10-05 14:35:15.049633","2021-10-05 14:35:15.049635","2021-10-05 14:35:15.049638","2021-10-05 14:35:15.049641","2021-10-05 14:35:15.049643","2021-10-05 14:35:15.049646","2021-10-05 14:35:15.049648","2021-10-05 14:35:15.049651","2021-10-05 14:35:15.049654","2021-10-05 14:35:15.049656","2021-10-05 14:35:15.049659","2021-10-05 14:35:15.049662","2021-10-05 14:35:15.049664","2021-10-05 14:35:15.049667","2021-10-05 14:35:15.049669","2021-10-05 14:35:15.049672","2021-10-05 14:35:15.049675","2021-10-05 14:35:15.049677","2021-10-05 14:35:15.049680","2021-10-05 14:35:15.049683","2021-10-05 14:35:15.049686","2021-10-05 14:35:15.049689","2021-10-05 14:35:15.049691","2021-10-05 14:35:15.049694","2021-10-05 14:35:15.049697","2021-10-05 14:35:15.049699","2021-10-05 14:35:15.049702","2021-10-05 14:35:15.049705","2021-10-05 14:35:15.049707","2021-10-05 14:35:15.049710","2021-10-05 14:35:15.049713","2021-10-05 14:35:15.049715","2021-10-05 14:35:15.049718","2021-10-05 14:35:15.049720","2021-10-05 14:35:15.049723","2021-10-05 14:35:15.049725","2021-10-05 14:35:15.049728")

for ($i = 0; $i -lt $rows.Length; $i++) {
    $dataSheet.Range("F" + $rows[$i]).Value = $timestamps[$i]
}

# --- 2. Add "metadata" sheet after "data" -----------------------------------------------------
$meta = $wb.Worksheets.Add($null, $dataSheet)
$meta.Name = "metadata"

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"
$meta.Range("A2").Value = 0

# Match header / index styling (bold, bordered, centered) used on the "data" sheet
$dataSheet.Range("B1:F1").Copy()
$meta.Range("B1:F1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$meta.Range("G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

$meta.Range("B2").Value = "Pierre Robin Sequence"
$meta.Range("C2").Value = 160
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "0.38"
$meta.Range("E2").Value = "2021-08-02T11:04:02.181399Z"
$meta.Range("F2").Value = "2021-10-05 14:35:15.045818"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/160/?format=json"

# Match page margins used elsewhere in the workbook (in points: 0.75in=54, 1in=72, 0.5in=36)
$meta.PageSetup.LeftMargin = 54
$meta.PageSetup.RightMargin = 54
$meta.PageSetup.TopMargin = 72
$meta.PageSetup.BottomMargin = 72
$meta.PageSetup.HeaderMargin = 36
$meta.PageSetup.FooterMargin = 36
